# MOSIP-14336: Updating Masterdata utility (language.xlsx)
# Reshapes the "language" sheet from columns A-E (code, name, family,
# nativeName, isActive) into a pandas-exported style table in columns
# B-G (Unnamed: 0, code, name, family, native_name, is_active), with an
# index column in A, and populates it with English/Arabic/French rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the bold/bordered/center-top header style that currently lives
# on A1 by copying its formatting onto the new header cells (B1:G1) and
# onto the new index column cells (A2:A4), before touching any values.
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# A1 is no longer used in the new layout - remove it (value + formatting).
$ws.Range("A1").Clear()
# The old data row (row 2, columns A:E) is being fully replaced.
$ws.Range("A2:E2").ClearContents()

# New header row
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "family"
$ws.Range("F1").Value = "native_name"
$ws.Range("G1").Value = "is_active"

# Row 2: English
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "eng"
$ws.Range("D2").Value = "English"
$ws.Range("E2").Value = "Indo-European"
$ws.Range("F2").Value = "English"
$ws.Range("G2").Value = $true

# Row 3: Arabic
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "ara"
$ws.Range("D3").Value = "Arabic"
$ws.Range("E3").Value = "Afro-Asiatic"
$ws.Range("F3").Value = "العَرَبِيَّة‎"
$ws.Range("G3").Value = $true

# Row 4: French
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "fra"
$ws.Range("D4").Value = "French"
$ws.Range("E4").Value = "Indo-European"
$ws.Range("F4").Value = "français"
$ws.Range("G4").Value = $true
